$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-45 (coin identity unchanged)
$ws.Range("D2").Value = "25.835.99"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.641.23"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "216.18"
$ws.Range("E5").Value = "  +0.59%  "

$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "0.0637"
$ws.Range("E9").Value = "  -0.62%  "

$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  -1.46%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.867.86"
$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("D14").Value = "1.642.42"
$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "0.563"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "63.11"
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").Value = "25.883.61"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("E19").Value = "  -0.15%  "

$ws.Range("D20").Value = "4.47"
$ws.Range("E20").Value = "  +2.39%  "

$ws.Range("D21").Value = "192.97"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "9.97"
$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  +3.09%  "

$ws.Range("E24").Value = "  +5.22%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").Value = "142.51"
$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  +2.12%  "

$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +0.83%  "

$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").Value = "3.25"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +0.94%  "

$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.910"
$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("D37").Value = "1.135.41"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("D39").Value = "0.547"
$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").Value = "0.0156"
$ws.Range("E40").Value = "  -0.35%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").Value = "5.63"
$ws.Range("E42").Value = "  +2.43%  "

$ws.Range("D43").Value = "100.75"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").Value = "0.806"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("D45").Value = "1.776.74"
$ws.Range("E45").Value = "  +0.32%  "

# Rows 46-51: coin list re-ranked (BabyDogeCoin dropped out of top list; others shift up; SynthetixNetwork newly enters)
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "55.45"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.419"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("E48").Value = "  +6.36%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +3.79%  "

$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  -2.06%  "
